$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue "D2" "24.699.64"
Set-TextValue "E2" "  -1.28%  "

Set-TextValue "D3" "1.678.47"
Set-TextValue "E3" "  -1.80%  "

Set-TextValue "D4" "1.003"
Set-TextValue "E4" "  +0.52%  "

Set-TextValue "D5" "313.80"
Set-TextValue "E5" "  -1.38%  "

Set-TextValue "D6" "1.001"
Set-TextValue "E6" "  +0.01%  "

Set-TextValue "D7" "0.3910"
Set-TextValue "E7" "  -3.25%  "

Set-TextValue "D8" "0.3954"
Set-TextValue "E8" "  -3.22%  "

Set-TextValue "B9" "BinanceUSD"
Set-TextValue "C9" "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
Set-TextValue "D9" "1.003"
Set-TextValue "E9" "  +0.36%  "

Set-TextValue "B10" "OKB"
Set-TextValue "C10" "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue "D10" "51.91"
Set-TextValue "E10" "  -3.30%  "

Set-TextValue "D11" "1.399"
Set-TextValue "E11" "  -5.68%  "

Set-TextValue "D12" "0.08633"
Set-TextValue "E12" "  -2.43%  "

Set-TextValue "D13" "25.24"
Set-TextValue "E13" "  -5.19%  "

Set-TextValue "D14" "7.315"
Set-TextValue "E14" "  -2.53%  "

Set-TextValue "B15" "ShibaInu"
Set-TextValue "C15" "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
Set-TextValue "D15" "0.00001316"
Set-TextValue "E15" "  -3.50%  "

Set-TextValue "B16" "Chainlink"
Set-TextValue "C16" "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-TextValue "D16" "7.757"
Set-TextValue "E16" "  -4.97%  "

Set-TextValue "D17" "1.679.67"
Set-TextValue "E17" "  -1.60%  "

Set-TextValue "D18" "93.63"
Set-TextValue "E18" "  -3.73%  "

Set-TextValue "D19" "0.07054"
Set-TextValue "E19" "  -1.61%  "

Set-TextValue "D20" "20.30"
Set-TextValue "E20" "  -4.56%  "

Set-TextValue "D21" "7.065"
Set-TextValue "E21" "  -3.05%  "

Set-TextValue "D22" "1.004"
Set-TextValue "E22" "  +0.30%  "

Set-TextValue "D23" "13.91"
Set-TextValue "E23" "  -3.49%  "

Set-TextValue "D24" "24.714.65"
Set-TextValue "E24" "  -1.14%  "

Set-TextValue "D25" "2.349"
Set-TextValue "E25" "  +0.84%  "

Set-TextValue "D26" "2.787"
Set-TextValue "E26" "  -5.10%  "

Set-TextValue "D27" "23.36"
Set-TextValue "E27" "  +0.05%  "

Set-TextValue "D28" "162.23"
Set-TextValue "E28" "  -2.87%  "

Set-TextValue "D29" "5.838"
Set-TextValue "E29" "  -7.25%  "

Set-TextValue "D30" "146.97"
Set-TextValue "E30" "  +0.88%  "

Set-TextValue "D31" "7.876"
Set-TextValue "E31" "  -6.87%  "

Set-TextValue "D32" "2.426"
Set-TextValue "E32" "  +8.52%  "

Set-TextValue "D33" "1.864.79"
Set-TextValue "E33" "  -1.77%  "

Set-TextValue "D34" "0.08406"
Set-TextValue "E34" "  -5.07%  "

Set-TextValue "D35" "0.03046"
Set-TextValue "E35" "  -4.92%  "

Set-TextValue "D36" "6.934"
Set-TextValue "E36" "  -4.41%  "

Set-TextValue "B37" "Algorand"
Set-TextValue "C37" "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextValue "D37" "0.2814"
Set-TextValue "E37" "  -2.30%  "

Set-TextValue "B38" "ImmutableX"
Set-TextValue "C38" "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue "D38" "0.9968"
Set-TextValue "E38" "  -3.81%  "

Set-TextValue "D39" "0.09478"
Set-TextValue "E39" "  +1.44%  "

Set-TextValue "D40" "10.58"
Set-TextValue "E40" "  -3.35%  "

Set-TextValue "D41" "1.505"
Set-TextValue "E41" "  +2.43%  "

Set-TextValue "D42" "0.7932"
Set-TextValue "E42" "  -6.76%  "

Set-TextValue "D43" "13.56"
Set-TextValue "E43" "  -4.64%  "

Set-TextValue "D44" "16.67"
Set-TextValue "E44" "  -4.41%  "

Set-TextValue "D45" "0.7136"
Set-TextValue "E45" "  -4.32%  "

Set-TextValue "D46" "2.562"
Set-TextValue "E46" "  -5.96%  "

Set-TextValue "D47" "4.193"
Set-TextValue "E47" "  -1.33%  "

Set-TextValue "D48" "0.08692"
Set-TextValue "E48" "  +3.90%  "

Set-TextValue "E49" "  +0.12%  "

Set-TextValue "D50" "1.346"
Set-TextValue "E50" "  -4.42%  "

Set-TextValue "D51" "138.01"
Set-TextValue "E51" "  -2.74%  "
